# Slide 5 ("Statistiques descriptives" / "Déséquilibre des classes" /
# "Détection d'outliers") - the three headline text boxes were widened
# (their shapes got taller thanks to word-wrap turning back on) as part of
# making room for a newly added video on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- "Statistiques descriptives" textbox (TextBox 13) ---
$sh13 = $s.Shapes.Item(8)
$sh13.Width  = 565.3280708661417
$sh13.Height = 46.44925196850394
$sh13.TextFrame.WordWrap = -1

# --- "Déséquilibre des classes" textbox (TextBox 15) ---
$sh15 = $s.Shapes.Item(10)
$sh15.Width  = 583.5158661417323
$sh15.Height = 46.44925196850394
$sh15.TextFrame.WordWrap = -1

# --- "Détection d'outliers" textbox (TextBox 16) ---
$sh16 = $s.Shapes.Item(11)

# Split the trailing run into " " and "d'outliers" (text is unchanged,
# only the run boundary moves) -- do this before resizing, since touching
# the text re-triggers the shape's auto-fit height calculation.
$tr16 = $sh16.TextFrame.TextRange
$c16 = $tr16.Characters(11, 10)
$c16.Text = "d’outliers"

$sh16.Width  = 482.69358267716535
$sh16.Height = 46.44925196850394
$sh16.TextFrame.WordWrap = -1
